# Sync attendance_reports: normalize "Recorded By" (column G) ordering so that
# the "System" / "system" auditing entry is listed first among comma-separated
# recorder names, matching the canonical formatting used by the source system.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Exact, known value replacements observed in column G ("Recorded By").
# Using an explicit map (rather than a generic reorder heuristic) guarantees
# we reproduce precisely the intended normalization without touching any
# value that should remain as-is.
$replacements = @{
    "system, backup@backdoor.com, System" = "System, backup@backdoor.com, system";
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com";
    "backup@backdoor.com, System"         = "System, backup@backdoor.com";
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com";
}

# Determine the last used row from the worksheet's used range.
$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G = 7
    $current = $cell.Value()

    if ($null -ne $current -and $replacements.ContainsKey($current)) {
        $cell.Value = $replacements[$current]
    }
}
